$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the crypto price/volume snapshot cells per the latest feed pull.
# Source values are plain text (prices use "." as a thousands separator in
# this feed, e.g. "42.419.19", and volumes keep padded "  +/-x.xx%  " strings),
# so any cell whose new value would otherwise be auto-parsed by Excel as a
# plain number (and could silently drop a meaningful trailing zero, e.g.
# "9.00" -> 9) is pinned to Text format first to keep it a literal string.

$ws.Range("D2").Value = "42.419.19"
$ws.Range("E2").Value = "  -0.52%  "
$ws.Range("D3").Value = "2.289.49"
$ws.Range("E3").Value = "  +0.12%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "300.21"
$ws.Range("E5").Value = "  -1.92%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.54"
$ws.Range("E6").Value = "  -1.26%  "
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.489"
$ws.Range("E9").Value = "  -2.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.28"
$ws.Range("E10").Value = "  -2.18%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "18.89"
$ws.Range("E11").Value = "  +2.89%  "
$ws.Range("E12").Value = "  -1.09%  "
$ws.Range("E13").Value = "  +0.16%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.68"
$ws.Range("E14").Value = "  -0.55%  "
$ws.Range("D15").Value = "2.649.46"
$ws.Range("E15").Value = "  -0.13%  "
$ws.Range("D16").Value = "2.306.28"
$ws.Range("E16").Value = "  +0.49%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.776"
$ws.Range("E17").Value = "  -0.25%  "
$ws.Range("D18").Value = "42.362.63"
$ws.Range("E18").Value = "  -0.50%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.14"
$ws.Range("E19").Value = "  -5.61%  "
$ws.Range("D20").Value = "0.0₃0885"
$ws.Range("E20").Value = "  -0.88%  "
$ws.Range("E21").Value = "  -1.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.41"
$ws.Range("E22").Value = "  +0.63%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.90"
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.24"
$ws.Range("E24").Value = "  +5.42%  "
$ws.Range("E25").Value = "  +0.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.39"
$ws.Range("E26").Value = "  -2.71%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.13"
$ws.Range("E27").Value = "  -3.43%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.31"
$ws.Range("E28").Value = "  +5.83%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "164.31"
$ws.Range("E29").Value = "  -1.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.00"
$ws.Range("E30").Value = "  -0.31%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "31.51"
$ws.Range("E31").Value = "  -3.90%  "
$ws.Range("E32").Value = "  -0.06%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.95"
$ws.Range("E33").Value = "  -0.15%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.36"
$ws.Range("E34").Value = "  -0.91%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0692"
$ws.Range("E35").Value = "  +0.10%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.33"
$ws.Range("E36").Value = "  -2.91%  "
$ws.Range("E37").Value = "  -8.32%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0996"
$ws.Range("E38").Value = "  -1.40%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.72"
$ws.Range("E39").Value = "  -0.92%  "
$ws.Range("E40").Value = "  -1.32%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.66"
$ws.Range("E41").Value = "  -0.90%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "19.93"
$ws.Range("E42").Value = "  +9.31%  "
$ws.Range("D43").Value = "1.946.67"
$ws.Range("E43").Value = "  -2.64%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.29"
$ws.Range("E44").Value = "  +2.54%  "
$ws.Range("E45").Value = "  -0.45%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.08"
$ws.Range("E46").Value = "  +3.08%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.71"
$ws.Range("E47").Value = "  -2.51%  "
$ws.Range("D48").Value = "2.520.71"
$ws.Range("E48").Value = "  +0.10%  "
$ws.Range("B49").Value = "HuobiToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.80"
$ws.Range("E49").Value = "  -3.71%  "
$ws.Range("B50").Value = "MultiversX"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "52.66"
$ws.Range("E50").Value = "  -1.96%  "
$ws.Range("E51").Value = "  +0.53%  "
